$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1878
$ws.Range("F6").Value = 726
$ws.Range("F9").Value = 875
$ws.Range("F10").Value = 1582
$ws.Range("F11").Value = 1265
$ws.Range("F12").Value = 1512
$ws.Range("F13").Value = 58
$ws.Range("F14").Value = 1476
$ws.Range("F18").Value = 1093
$ws.Range("F19").Value = 362
$ws.Range("F22").Value = 1695
$ws.Range("F23").Value = 214
$ws.Range("F26").Value = 1198
$ws.Range("F27").Value = 311558
$ws.Range("F32").Value = 1141
$ws.Range("F35").Value = 1130
$ws.Range("F36").Value = 1081
$ws.Range("F37").Value = 258
$ws.Range("F38").Value = 70
$ws.Range("F42").Value = 116
$ws.Range("F43").Value = 2018
$ws.Range("F44").Value = 87
$ws.Range("F47").Value = 796

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 183
$ws.Range("F9").Value = 2584
$ws.Range("F18").Value = 459
$ws.Range("F19").Value = 23
$ws.Range("F21").Value = 313
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 26
$ws.Range("F24").Value = 3
$ws.Range("F30").Value = 213
$ws.Range("F31").Value = 60
$ws.Range("F33").Value = 57
$ws.Range("F39").Value = 60
$ws.Range("F40").Value = 60

$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 2880
$ws.Range("F6").Value = 4640
$ws.Range("F7").Value = 136
$ws.Range("F10").Value = 732
$ws.Range("F11").Value = 467
$ws.Range("F12").Value = 349
$ws.Range("F13").Value = 1071
$ws.Range("F14").Value = 283
$ws.Range("F15").Value = 665

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1878
$ws.Range("F4").Value = 4640
$ws.Range("F5").Value = 732
$ws.Range("F6").Value = 467
$ws.Range("F7").Value = 349
$ws.Range("F8").Value = 349
$ws.Range("F9").Value = 1071
$ws.Range("F11").Value = 875
$ws.Range("F13").Value = 1582
$ws.Range("F14").Value = 1265
$ws.Range("F15").Value = 1512
$ws.Range("F16").Value = 1477
$ws.Range("F20").Value = 1093
$ws.Range("F21").Value = 362
$ws.Range("F22").Value = 665
$ws.Range("F23").Value = 666
$ws.Range("F24").Value = 459
$ws.Range("F25").Value = 1695
$ws.Range("F26").Value = 214
$ws.Range("F29").Value = 1198
$ws.Range("F30").Value = 313
$ws.Range("F33").Value = 1141
$ws.Range("F36").Value = 1130
$ws.Range("F37").Value = 3
$ws.Range("F38").Value = 1081
$ws.Range("F39").Value = 259
$ws.Range("F40").Value = 874
$ws.Range("F41").Value = 60
$ws.Range("F44").Value = 116
$ws.Range("F45").Value = 2018
$ws.Range("F46").Value = 87
$ws.Range("F48").Value = 60
$ws.Range("F49").Value = 796
